$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New test case row (LOG-TC-05 / "Login with empty password") ---
# Copy column formatting D:L from the row above (row 7) down into row 8
# so the new row picks up the same per-column styles (s=2/2/2/2/2/3/4/2/2).
$ws.Range("D7:L7").Copy()
$ws.Range("D8:L8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D8").Value2 = "LOG-TC-05"
$ws.Range("E8").Value2 = "Login Module"
$ws.Range("F8").Value2 = "LOG-TS-04"
$ws.Range("G8").Value2 = "Login with empty password"
$ws.Range("H8").Value2 = "On login page "
$ws.Range("I8").Value2 = "Enter email, leave password empty, click Login"
$ws.Range("J8").Value2 = "user1@gmail.com" + [char]10 + '""'
$ws.Range("K8").Value2 = "Validation message appears"
$ws.Range("L8").Value2 = "High"

$ws.Rows(8).RowHeight = 47.25

# Hyperlink for the new test-data cell (J8), mirroring the existing
# mailto-style hyperlinks used for the other test-data cells.
$ws.Hyperlinks.Add($ws.Range("J8"), 'mailto:user1@gmail.com%0a""')

# --- Merge the SUB / Scenario cells of the two related rows (7 & 8) ---
# Merging adjusts B7/C7 to the "top of merge" border style and creates the
# matching "bottom of merge" style on B8/C8 automatically.
$ws.Range("B7:B8").Merge()
$ws.Range("C7:C8").Merge()

# --- Column D width tweak (drop auto bestFit, use an explicit width) ---
$ws.Columns("D").ColumnWidth = 12.74

# --- Sheet view: scroll / selection change ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("K11").Select()
